$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume table update (GitHub Actions data refresh)

$ws.Range("D2").Value = "42.089.69"
$ws.Range("E2").Value = "  +2.04%  "
$ws.Range("D3").Value = "2.219.17"
$ws.Range("E3").Value = "  +1.25%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.616"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.76%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "68.04"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.90%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.620"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.42"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.33"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.11%  "
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.91%  "
$ws.Range("E14").Value = "  -0.38%  "
$ws.Range("D15").Value = "2.555.16"
$ws.Range("E15").Value = "  +1.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.871"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("D18").Value = "2.227.20"
$ws.Range("E18").Value = "  +1.98%  "
$ws.Range("D19").Value = "42.033.80"
$ws.Range("E19").Value = "  +1.80%  "
$ws.Range("D20").Value = "0.0₃0960"
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("E22").Value = "  -1.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.17%  "
$ws.Range("E25").Value = "  -0.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.41"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.62%  "
$ws.Range("E29").Value = "  -1.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "167.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.46"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.15"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +10.90%  "
$ws.Range("E34").Value = "  +1.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0781"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.37%  "
$ws.Range("E36").Value = "  -0.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.61"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0316"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.23%  "
$ws.Range("E41").Value = "  +1.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.85"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.35%  "
$ws.Range("E46").Value = "  -3.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("E48").Value = "  -0.87%  "
$ws.Range("E49").Value = "  -0.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.30%  "
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.07%  "
